$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 3437  # was 3441
$ws.Cells.Item(3, 4).Value = 3931  # was 3935
$ws.Cells.Item(4, 4).Value = 3548  # was 3552
$ws.Cells.Item(5, 4).Value = 3527  # was 3531
$ws.Cells.Item(6, 4).Value = 3809  # was 3814
$ws.Cells.Item(7, 4).Value = 3946  # was 3950
$ws.Cells.Item(8, 4).Value = 4176  # was 4181
$ws.Cells.Item(9, 4).Value = 4483  # was 4488
$ws.Cells.Item(10, 4).Value = 4895  # was 4900
$ws.Cells.Item(11, 4).Value = 5533  # was 5539
$ws.Cells.Item(12, 4).Value = 5854  # was 5860
$ws.Cells.Item(13, 4).Value = 6400  # was 6407
$ws.Cells.Item(14, 4).Value = 7108  # was 7116
$ws.Cells.Item(15, 4).Value = 7643  # was 7652
$ws.Cells.Item(16, 4).Value = 8077  # was 8086
$ws.Cells.Item(17, 4).Value = 8845  # was 8855
$ws.Cells.Item(18, 4).Value = 9492  # was 9502
$ws.Cells.Item(19, 4).Value = 10233  # was 10244
$ws.Cells.Item(20, 4).Value = 10660  # was 10672
$ws.Cells.Item(21, 4).Value = 10628  # was 10643
$ws.Cells.Item(22, 4).Value = 11314  # was 11326
$ws.Cells.Item(23, 4).Value = 11814  # was 11821
$ws.Cells.Item(24, 4).Value = 12235  # was 12244
$ws.Cells.Item(25, 4).Value = 12851  # was 12846
$ws.Cells.Item(26, 4).Value = 14001  # was 13997
$ws.Cells.Item(27, 4).Value = 15117  # was 15109
$ws.Cells.Item(28, 4).Value = 16401  # was 16383
$ws.Cells.Item(29, 4).Value = 17489  # was 17468
$ws.Cells.Item(30, 4).Value = 18262  # was 18244
$ws.Cells.Item(31, 4).Value = 17890  # was 17895
$ws.Cells.Item(32, 4).Value = 18954  # was 18952
$ws.Cells.Item(35, 4).Value = 22387  # was 22386
$ws.Cells.Item(38, 4).Value = 22752  # was 22751
$ws.Cells.Item(39, 4).Value = 23720  # was 23718
$ws.Cells.Item(41, 4).Value = 24968  # was 24969
$ws.Cells.Item(42, 4).Value = 23363  # was 23366
$ws.Cells.Item(43, 4).Value = 26513  # was 24928
$ws.Cells.Item(44, 4).Value = 27629  # was 26195
$ws.Cells.Item(45, 4).Value = 28546  # was 27246
$ws.Cells.Item(46, 4).Value = 29507  # was 28252
$ws.Cells.Item(47, 4).Value = 30554  # was 29248
$ws.Cells.Item(48, 4).Value = 31674  # was 30279
